# "added img path added verify err msg method created invalid login script"
#
# The login fixture's "testValidLogin" sheet gets a deliberately-wrong
# password in B2 ("manager" -> "damager") so it can double as the data
# for the newly added invalid-login test script, and the sheet's saved
# selection moves from A2 to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testValidLogin")
$ws.Activate()

$ws.Range("B2").Value = "damager"

$null = $ws.Range("B3").Select()
